# Add two new columns, I (I0) and J (IF), to the sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row values
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy the header formatting (bold, centered, bordered) from H1 onto the
# new header cells so they match the rest of the header row.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# Data rows 2..23 for columns I and J
$values = @(
    @(8, 8),
    @(5, 6),
    @(10, 10),
    @(8, 8),
    @(7, 8),
    @(8, 8),
    @(8, 8),
    @(5, 5),
    @(6, 6),
    @(5, 6),
    @(5, 5),
    @(7, 7),
    @(6, 7),
    @(4, 4),
    @(4, 5),
    @(5, 5),
    @(7, 7),
    @(8, 8),
    @(7, 7),
    @(8, 8),
    @(6, 6),
    @(4, 4)
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $values[$i][0]
    $ws.Cells.Item($row, 10).Value = $values[$i][1]
}
